$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Parneet Kaur"

# Row 7 - __init__ / Valid input values
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = '(1010, "Susan", "Clark", "susan@example.com")'
$ws.Range("G7").Value = "Object created successfully. All private attributes set to correct values."

# Row 8 - __init__ / Invalid client number
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = '("ABC", "Susan", "Clark", "susan@example.com")'
$ws.Range("G8").Value = 'ValueError raised: "Client number must be an integer."'

# Row 9 - __init__ / Blank first name
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = '(1011, " ", "Clark", "susan@example.com")'
$ws.Range("G9").Value = 'ValueError raised: "First name cannot be blank."'

# Row 10 - __init__ / Blank last name
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = '(1012, "Susan", " ", "susan@example.com")'
$ws.Range("G10").Value = 'ValueError raised: "Last name cannot be blank."'

# Row 11 - __init__ / Invalid email
$ws.Range("E11").Value = "None"
$ws.Range("F11").Value = '(2020, "John", "Smith", "invalidemail")'
$ws.Range("G11").Value = 'Email address set to "email@pixell-river.com"'

# Row 12 - client_number getter
$ws.Range("E12").Value = "Valid Client object exists"
$ws.Range("F12").Value = "obj.client_number"
$ws.Range("G12").Value = "Returns integer client number (e.g. 1010)."

# Row 13 - first_name getter
$ws.Range("E13").Value = "Valid Client object exists"
$ws.Range("F13").Value = "obj.first_name"
$ws.Range("G13").Value = 'Returns string "Susan".'

# Row 14 - last_name getter
$ws.Range("E14").Value = "Valid Client object exists"
$ws.Range("F14").Value = "obj.last_name"
$ws.Range("G14").Value = 'Returns string "Clark".'

# Row 15 - email_address getter
$ws.Range("E15").Value = "Valid Client object exists"
$ws.Range("F15").Value = "obj.email_address"
$ws.Range("G15").Value = 'Returns valid email string "susan@example.com".'

# Row 16 - __str__
$ws.Range("E16").Value = "Valid Client object exists"
$ws.Range("F16").Value = "str(obj)"
$ws.Range("G16").Value = '"Clark, Susan [1010] - susan@example.com" followed by newline.'

# View changes: zoom to 87% and select L23
$ws.Application.ActiveWindow.Zoom = 87
[void]$ws.Range("L23").Select()
